# Update data from 2024 and 2025 and some minor changes to the codebook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New party rows appended to the bottom of the table (rows 23-31)
$ws.Range("A23").Value = "BD"
$ws.Range("B23").Value = "Bündnis Deutschland"
$ws.Range("H23").Value = "https://de.wikipedia.org/wiki/B%C3%BCndnis_Deutschland"

$ws.Range("A24").Value = "BSW"
$ws.Range("B24").Value = "Bündnis Sahra Wagenknecht"
$ws.Range("D24").Value = 313
$ws.Range("H24").Value = "https://de.wikipedia.org/wiki/B%C3%BCndnis_Sahra_Wagenknecht"

$ws.Range("A25").Value = "FS (2021)"
$ws.Range("B25").Value = "Freie Sachsen (2021)"
$ws.Range("G25").Value = "Not to be confused with the previous party also called Freie Sachsen."
$ws.Range("H25").Value = "https://de.wikipedia.org/wiki/Freie_Sachsen_(2021)"

$ws.Range("A26").Value = "WU"
$ws.Range("B26").Value = "WerteUnion"
$ws.Range("H26").Value = "https://de.wikipedia.org/wiki/Werteunion"

$ws.Range("A27").Value = "PB"
$ws.Range("B27").Value = "Plus Brandenburg"
$ws.Range("G27").Value = "Joint list of Pirate Party, ÖDP and Volt for 2024 Brandenburg election"
$ws.Range("H27").Value = "https://de.wikipedia.org/wiki/Plus_Brandenburg"

$ws.Range("A28").Value = "DLW"
$ws.Range("B28").Value = "Deutsch Land Wirtschaft"
$ws.Range("H28").Value = "https://de.wikipedia.org/wiki/Deutsch_Land_Wirtschaft"

$ws.Range("B29").Value = "Die Wahl für Frieden und soziale Gerechtigkeit"
$ws.Range("A29").Value = "WFG"
$ws.Range("H29").Value = "https://de.wikipedia.org/wiki/Die_Wahl_f%C3%BCr_Frieden_und_soziale_Gerechtigkeit"

$ws.Range("A30").Value = "DAVA"
$ws.Range("B30").Value = "Demokratische Allianz für Vielfalt und Aufbruch"
$ws.Range("H30").Value = "https://de.wikipedia.org/wiki/Demokratische_Allianz_f%C3%BCr_Vielfalt_und_Aufbruch"

$ws.Range("A31").Value = "NPD (2023)"
$ws.Range("B31").Value = "Nationaldemokratische Partei Deutschlands (2023)"
$ws.Range("H31").Value = "https://de.wikipedia.org/wiki/Nationaldemokratische_Partei_Deutschlands_(2023)"
$ws.Range("G31").Value = "Split from the old NPD after the party renamed itself to Die Heimat"

# Scroll/selection state like in the authored file
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("E44").Select()
